$d = $word.ActiveDocument

# Replace the placeholder text with the actual URL
$d.Content.Find.Execute("(Mettre le lien URL ici)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "https://apepingagne.github.io/420-KB2-PFI-Alexandre-Pepin-Gagn-/", 2)
